$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from SCD0286 to SCD0018
$ws.Name = "SCD0018"

# Update the TC_ID column (B) for rows 3-7 to the new test-case id.
# Row 2 keeps its original "DGS-301" value (untouched by this change).
$ws.Range("B3").Value = "SCD0018-009"
$ws.Range("B4").Value = "SCD0018-009"
$ws.Range("B5").Value = "SCD0018-009"
$ws.Range("B6").Value = "SCD0018-009"
$ws.Range("B7").Value = "SCD0018-009"

# Widen column B to fit the longer TC_ID text (target stored width ~12.5703125)
$ws.Columns.Item(2).ColumnWidth = 11.6

# Row 5 wraps onto fewer lines now that column B is wider
$ws.Rows.Item(5).RowHeight = 89.25

# Move/restore the active selection to B8
$ws.Range("B8").Select()
